# Update the EPEX Spot prices workbook with the latest data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": add a new day column AZ (04-aug) with hourly prices
# ---------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting (bold header, border, centered alignment) from the
# previous day's header cell onto the new one, then set its value.
$wsSpot.Range("AY1").Copy($wsSpot.Range("AZ1"))
$wsSpot.Range("AZ1").Value = "04-aug"

$spotValues = @(
    80.98999999999999,
    70.06999999999999,
    69.97,
    63.8,
    61.78,
    65.7,
    74.33,
    78.83,
    78.43000000000001,
    52.57,
    30,
    13.75,
    6.76,
    4.05,
    0.65,
    0,
    5.79,
    7.28,
    38.94,
    63.7,
    70.25,
    74,
    67.14,
    49.4
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 52).Value = $spotValues[$i]
}

# ---------------------------------------------------------------
# Sheet "Gaz": append two new daily rows (2025-08-02, 2025-08-03)
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date-like text to stay literal text (not auto-converted to a
# date serial) while keeping the cells on the default "Normal" style, same
# as the rest of the column.
$wsGaz.Range("A49:A50").NumberFormat = "@"
$wsGaz.Range("A49").Value = "2025-08-02"
$wsGaz.Range("B49").Value = 32.775
$wsGaz.Range("A50").Value = "2025-08-03"
$wsGaz.Range("B50").Value = 32.775
$wsGaz.Range("A49:A50").Style = "Normal"

# ---------------------------------------------------------------
# Sheet "CO2": append two new daily rows (2025-08-02, 2025-08-03)
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A49:A50").NumberFormat = "@"
$wsCo2.Range("A49").Value = "2025-08-02"
$wsCo2.Range("B49").Value = 70.58
$wsCo2.Range("A50").Value = "2025-08-03"
$wsCo2.Range("B50").Value = 70.58
$wsCo2.Range("A49:A50").Style = "Normal"
